# Fixing the big mistake: update the average/analysis statistics that were
# recomputed, along with the dependent Total/% energy sector figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# mean / std / min / 25% / 50% / 75% / max rows for "Total" (B) and "Community" (D)
$ws.Range("B3").Value = 4041.924321666666
$ws.Range("D3").Value = 271.2254994101976

$ws.Range("B4").Value = 1693.339467015477
$ws.Range("D4").Value = 93.9815526326003

$ws.Range("B5").Value = 1390.186172602739
$ws.Range("D5").Value = 94.74842465753429

$ws.Range("B6").Value = 2769.910641780824
$ws.Range("D6").Value = 185.0059212328762

$ws.Range("B7").Value = 3451.0174630137
$ws.Range("D7").Value = 273.8044246575339

$ws.Range("B8").Value = 5580.952449315071
$ws.Range("D8").Value = 357.5890547945205

$ws.Range("B9").Value = 7444.770926027398
$ws.Range("D9").Value = 432.2652630136973

# Total / Residential / Community / IGA sums and % energy sector
$ws.Range("F10").Value = 5820371.023200002

$ws.Range("G11").Value = 0.8063540400844781

$ws.Range("F12").Value = 390564.7191506845
$ws.Range("G12").Value = 0.06710306226078945

$ws.Range("G13").Value = 0.1265428976547325
